$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 118 - TRADING_ATTEMPT
$ws.Cells.Item(118, 1).Value = "2025-11-24T01:46:53.132304"
$ws.Cells.Item(118, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(118, 3).Value = "BTC"
$ws.Cells.Item(118, 4).Value = "UNKNOWN"
$ws.Cells.Item(118, 5).Value = 86576.78246280954
$ws.Cells.Item(118, 11).Value = "ATTEMPT"
$ws.Cells.Item(118, 12).Value = "Attempting trade 1/1"

# Row 119 - POSITION_OPENED
$ws.Cells.Item(119, 1).Value = "2025-11-24T01:46:54.474317"
$ws.Cells.Item(119, 2).Value = "POSITION_OPENED"
$ws.Cells.Item(119, 3).Value = "BTC"
$ws.Cells.Item(119, 4).Value = "UNKNOWN"
$ws.Cells.Item(119, 5).Value = 86576.78246280954
$ws.Cells.Item(119, 6).Value = 3600
$ws.Cells.Item(119, 7).Value = 40
$ws.Cells.Item(119, 8).Value = 0.3108185808726259
$ws.Cells.Item(119, 11).Value = "SUCCESS"
